$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 499.75
$ws.Range("I32").Value = 401
$ws.Range("J32").Value = 532.6667
$ws.Range("K32").Value = 401
$ws.Range("L32").Value = 532.6667
$ws.Range("M32").Value = -75
$ws.Range("N32").Value = -1184.6667
$ws.Range("H64").Value = 3959.16
$ws.Range("I64").Value = 3478.6
$ws.Range("K64").Value = 3478.6
$ws.Range("M64").Value = -3230.6
$ws.Range("H67").Value = 3959.16
$ws.Range("I67").Value = 3478.6
$ws.Range("K67").Value = 3478.6
$ws.Range("M67").Value = -2620.6
$ws.Range("H76").Value = 4632679.5
$ws.Range("I76").Value = 3228.5715
$ws.Range("J76").Value = 11113911
$ws.Range("K76").Value = 3228.5715
$ws.Range("L76").Value = 11113911
$ws.Range("M76").Value = -2913.5715
$ws.Range("N76").Value = -11114541
$ws.Range("H79").Value = 4632679.5
$ws.Range("I79").Value = 3228.5715
$ws.Range("J79").Value = 11113911
$ws.Range("K79").Value = 3228.5715
$ws.Range("L79").Value = 11113911
$ws.Range("M79").Value = -2136.5715
$ws.Range("N79").Value = -11116095
$ws.Range("H86").Value = 7703.1113
$ws.Range("I86").Value = 1998.75
$ws.Range("J86").Value = 12266.6
$ws.Range("K86").Value = 1998.75
$ws.Range("L86").Value = 12266.6
$ws.Range("M86").Value = -875.75
$ws.Range("N86").Value = -14512.6
$ws.Range("H89").Value = 7703.1113
$ws.Range("I89").Value = 1998.75
$ws.Range("J89").Value = 12266.6
$ws.Range("K89").Value = 9993.75
$ws.Range("L89").Value = 61333
$ws.Range("M89").Value = -4377.75
$ws.Range("N89").Value = -72565
$ws.Range("H135").Value = 20840886
$ws.Range("I135").Value = 790.6111
$ws.Range("J135").Value = 83361170
$ws.Range("K135").Value = 7115.4999
$ws.Range("L135").Value = 750250530
$ws.Range("M135").Value = -4580.4999
$ws.Range("N135").Value = -750255600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2389.739
$ws.Range("I61").Value = 1665.7778
$ws.Range("K61").Value = 1665.7778
$ws.Range("M61").Value = -1453.7778
$ws.Range("H63").Value = 3127150
$ws.Range("I63").Value = 2437.5
$ws.Range("J63").Value = 15626000
$ws.Range("K63").Value = 2437.5
$ws.Range("L63").Value = 15626000
$ws.Range("M63").Value = -1751.5
$ws.Range("N63").Value = -15627372
$ws.Range("H66").Value = 3127150
$ws.Range("I66").Value = 2437.5
$ws.Range("J66").Value = 15626000
$ws.Range("K66").Value = 12187.5
$ws.Range("L66").Value = 78130000
$ws.Range("M66").Value = -8755.5
$ws.Range("N66").Value = -78136864
$ws.Range("H132").Value = 25322.045
$ws.Range("I132").Value = 2503.6875
$ws.Range("J132").Value = 86171
$ws.Range("K132").Value = 7511.0625
$ws.Range("L132").Value = 258513
$ws.Range("M132").Value = -4981.0625
$ws.Range("N132").Value = -263573
$ws.Range("H136").Value = 2389.739
$ws.Range("I136").Value = 1665.7778
$ws.Range("K136").Value = 4997.3334
$ws.Range("M136").Value = -2447.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2438.4
$ws.Range("I94").Value = 1160.8889
$ws.Range("K94").Value = 1160.8889
$ws.Range("M94").Value = -709.8888999999999
$ws.Range("H134").Value = 29711.316
$ws.Range("I134").Value = 41282.594
$ws.Range("J134").Value = 1309.091
$ws.Range("K134").Value = 123847.782
$ws.Range("L134").Value = 3927.273
$ws.Range("M134").Value = -121312.782
$ws.Range("N134").Value = -8997.272999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19761.777
$ws.Range("I31").Value = 30710.273
$ws.Range("J31").Value = 2557
$ws.Range("K31").Value = 30710.273
$ws.Range("L31").Value = 2557
$ws.Range("M31").Value = -30415.273
$ws.Range("N31").Value = -3147
$ws.Range("H34").Value = 19761.777
$ws.Range("I34").Value = 30710.273
$ws.Range("J34").Value = 2557
$ws.Range("K34").Value = 30710.273
$ws.Range("L34").Value = 2557
$ws.Range("M34").Value = -30508.273
$ws.Range("N34").Value = -2961
$ws.Range("H62").Value = 142862060
$ws.Range("I62").Value = 333336000
$ws.Range("J62").Value = 6600
$ws.Range("K62").Value = 333336000
$ws.Range("L62").Value = 6600
$ws.Range("M62").Value = -333335376
$ws.Range("N62").Value = -7848
$ws.Range("H65").Value = 142862060
$ws.Range("I65").Value = 333336000
$ws.Range("J65").Value = 6600
$ws.Range("K65").Value = 1666680000
$ws.Range("L65").Value = 33000
$ws.Range("M65").Value = -1666676880
$ws.Range("N65").Value = -39240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 750
$ws.Range("I98").Value = 600
$ws.Range("J98").Value = 775
$ws.Range("K98").Value = 1800
$ws.Range("L98").Value = 2325
$ws.Range("M98").Value = -302
$ws.Range("N98").Value = -5321
$ws.Range("H131").Value = 741.51514
$ws.Range("J131").Value = 741.51514
$ws.Range("L131").Value = 2224.54542
$ws.Range("N131").Value = -12304.54542
$ws.Range("H139").Value = 2805.3076
$ws.Range("I139").Value = 1853.625
$ws.Range("J139").Value = 4328
$ws.Range("K139").Value = 5560.875
$ws.Range("L139").Value = 12984
$ws.Range("M139").Value = -420.875
$ws.Range("N139").Value = -23264
$ws.Range("H140").Value = 1591
$ws.Range("I140").Value = 1333.0769
$ws.Range("K140").Value = 3999.2307
$ws.Range("M140").Value = 1180.7693
$ws.Range("H141").Value = 3464.8333
$ws.Range("I141").Value = 5517
$ws.Range("K141").Value = 16551
$ws.Range("M141").Value = -11371

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6959900.5
$ws.Range("I70").Value = 38333.332
$ws.Range("J70").Value = 10420684
$ws.Range("K70").Value = 38333.332
$ws.Range("L70").Value = 10420684
$ws.Range("M70").Value = -38063.332
$ws.Range("N70").Value = -10421224
$ws.Range("H73").Value = 6959900.5
$ws.Range("I73").Value = 38333.332
$ws.Range("J73").Value = 10420684
$ws.Range("K73").Value = 38333.332
$ws.Range("L73").Value = 10420684
$ws.Range("M73").Value = -37397.332
$ws.Range("N73").Value = -10422556
$ws.Range("H80").Value = 3785.4285
$ws.Range("I80").Value = 3449.6667
$ws.Range("J80").Value = 4037.25
$ws.Range("K80").Value = 3449.6667
$ws.Range("L80").Value = 4037.25
$ws.Range("M80").Value = -2451.6667
$ws.Range("N80").Value = -6033.25
$ws.Range("H83").Value = 3785.4285
$ws.Range("I83").Value = 3449.6667
$ws.Range("J83").Value = 4037.25
$ws.Range("K83").Value = 17248.3335
$ws.Range("L83").Value = 20186.25
$ws.Range("M83").Value = -12256.3335
$ws.Range("N83").Value = -30170.25
$ws.Range("H113").Value = 2117.1936
$ws.Range("J113").Value = 2240.7896
$ws.Range("L113").Value = 2240.7896
$ws.Range("N113").Value = -6580.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 13722
$ws.Range("I74").Value = 15444
$ws.Range("J74").Value = 12000
$ws.Range("K74").Value = 15444
$ws.Range("L74").Value = 12000
$ws.Range("M74").Value = -14446
$ws.Range("N74").Value = -13996
$ws.Range("H77").Value = 13722
$ws.Range("I77").Value = 15444
$ws.Range("J77").Value = 12000
$ws.Range("K77").Value = 46332
$ws.Range("L77").Value = 36000
$ws.Range("M77").Value = -41340
$ws.Range("N77").Value = -45984
$ws.Range("H100").Value = 3164.1428
$ws.Range("I100").Value = 1722.2222
$ws.Range("K100").Value = 1722.2222
$ws.Range("M100").Value = -1181.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 250001230
$ws.Range("I81").Value = 1633.3334
$ws.Range("J81").Value = 1000000000
$ws.Range("K81").Value = 3266.6668
$ws.Range("L81").Value = 2000000000
$ws.Range("M81").Value = -2205.6668
$ws.Range("N81").Value = -2000002122
$ws.Range("H84").Value = 250001230
$ws.Range("I84").Value = 1633.3334
$ws.Range("J84").Value = 1000000000
$ws.Range("K84").Value = 16333.334
$ws.Range("L84").Value = 10000000000
$ws.Range("M84").Value = -11029.334
$ws.Range("N84").Value = -10000010608
